$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -0.1510908837489806
$ws.Cells.Item(2, 3).Value = 2.084195974016096
$ws.Cells.Item(2, 4).Value = 20.51875844210171
$ws.Cells.Item(2, 5).Value = 4.529763618788701
$ws.Cells.Item(2, 6).Value = 4.628991593011878
$ws.Cells.Item(2, 7).Value = 23

$ws.Cells.Item(3, 2).Value = -0.2243426211202774
$ws.Cells.Item(3, 3).Value = 2.382314769375817
$ws.Cells.Item(3, 4).Value = 19.09061844583581
$ws.Cells.Item(3, 5).Value = 4.369281227597487
$ws.Cells.Item(3, 6).Value = 4.46620300197608
$ws.Cells.Item(3, 7).Value = 22

$ws.Cells.Item(4, 2).Value = -0.7120868115500602
$ws.Cells.Item(4, 3).Value = 1.933618923903567
$ws.Cells.Item(4, 4).Value = 10.88254014902561
$ws.Cells.Item(4, 5).Value = 3.298869525917267
$ws.Cells.Item(4, 6).Value = 3.300643293046703
$ws.Cells.Item(4, 7).Value = 21

$ws.Cells.Item(5, 2).Value = -0.3051420100610381
$ws.Cells.Item(5, 3).Value = 1.578682252322016
$ws.Cells.Item(5, 4).Value = 10.2251943824045
$ws.Cells.Item(5, 5).Value = 3.197685785439917
$ws.Cells.Item(5, 6).Value = 3.265784782946779
$ws.Cells.Item(5, 7).Value = 20

$ws.Cells.Item(6, 2).Value = -0.2804117665781
$ws.Cells.Item(6, 3).Value = 1.683207935177004
$ws.Cells.Item(6, 4).Value = 11.04924039551323
$ws.Cells.Item(6, 5).Value = 3.324039770446983
$ws.Cells.Item(6, 6).Value = 3.402952828034286
$ws.Cells.Item(6, 7).Value = 19

$ws.Cells.Item(7, 2).Value = -0.3173356620072107
$ws.Cells.Item(7, 3).Value = 1.884213172660337
$ws.Cells.Item(7, 4).Value = 10.94775440033029
$ws.Cells.Item(7, 5).Value = 3.308739095234057
$ws.Cells.Item(7, 6).Value = 3.388969517186058
$ws.Cells.Item(7, 7).Value = 18

$ws.Cells.Item(8, 2).Value = -0.2356007239738259
$ws.Cells.Item(8, 3).Value = 1.940837229566887
$ws.Cells.Item(8, 4).Value = 12.01945225291491
$ws.Cells.Item(8, 5).Value = 3.466908169091721
$ws.Cells.Item(8, 6).Value = 3.565345857874665
$ws.Cells.Item(8, 7).Value = 17

$ws.Cells.Item(9, 2).Value = -0.1538742921933148
$ws.Cells.Item(9, 3).Value = 1.930586848953809
$ws.Cells.Item(9, 4).Value = 10.98523007623265
$ws.Cells.Item(9, 5).Value = 3.314397392623983
$ws.Cells.Item(9, 6).Value = 3.41940389010477
$ws.Cells.Item(9, 7).Value = 16

$ws.Cells.Item(10, 2).Value = -0.1230455590706185
$ws.Cells.Item(10, 3).Value = 1.907399787585887
$ws.Cells.Item(10, 4).Value = 12.46868919605704
$ws.Cells.Item(10, 5).Value = 3.531103113200893
$ws.Cells.Item(10, 6).Value = 3.652819212576486
$ws.Cells.Item(10, 7).Value = 15

$ws.Cells.Item(11, 2).Value = -0.09178185018074302
$ws.Cells.Item(11, 3).Value = 2.13619364861378
$ws.Cells.Item(11, 4).Value = 13.88641155430669
$ws.Cells.Item(11, 5).Value = 3.726447578365579
$ws.Cells.Item(11, 6).Value = 3.865944277603689
$ws.Cells.Item(11, 7).Value = 14

